# "checkbox functions implemented (110, 112)"
# Add a new "testbench" column (K) to the first worksheet (Tabelle1),
# assigning alternating SYS-110.tbc / SYS-112.tbc values to the 18 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell for the new column, right-aligned / numeric style like the
# other header cells in row 1 (matches the look of column J's header).
$ws.Range("K1").Value = "testbench"
$ws.Range("K1").NumberFormat = $ws.Range("J1").NumberFormat
$ws.Range("K1").HorizontalAlignment = $ws.Range("J1").HorizontalAlignment

# Data rows 2..19 alternate between SYS-110.tbc (even rows) and SYS-112.tbc (odd rows)
$values = @(
    "SYS-110.tbc","SYS-112.tbc","SYS-112.tbc","SYS-110.tbc","SYS-112.tbc",
    "SYS-112.tbc","SYS-110.tbc","SYS-112.tbc","SYS-110.tbc","SYS-112.tbc",
    "SYS-110.tbc","SYS-112.tbc","SYS-110.tbc","SYS-112.tbc","SYS-110.tbc",
    "SYS-112.tbc","SYS-110.tbc","SYS-112.tbc"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}

# Size the new column (closest achievable width to the authored 15.140625)
$ws.Columns.Item(11).ColumnWidth = 14.25

# Update selection to match the committed workbook state
[void]$ws.Range("K1:K19").Select()
